# Geelong_stats.xlsx — insert 4 duplicate "round" columns right after the
# current last data column (KV), matching the upcoming-predictions layout
# that a later guard-check will rely on.
#
# Before: data runs A1:KV102, with KV being the last (most-recent) round
#         column for every stat row.
# After:  4 new columns (KW:KZ) are inserted in place of the old KV, each
#         a duplicate of the values that used to live in KV. The old KV
#         column's data is pushed out to KZ by the insert, and the new
#         KV:KY columns are filled with the same values so every one of
#         the 5 right-most columns (KV,KW,KX,KY,KZ) ends up holding the
#         value that used to be in the single KV column.

function Get-ColumnLetter([int]$col) {
    $letter = ""
    while ($col -gt 0) {
        $rem = ($col - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $col = [int](($col - $rem - 1) / 26)
    }
    return $letter
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

$numNewCols = 4

$lastColLetter      = Get-ColumnLetter $lastCol
$insertEndLetter    = Get-ColumnLetter ($lastCol + $numNewCols - 1)
$lastNewColLetter   = Get-ColumnLetter ($lastCol + $numNewCols)

# Insert 4 blank columns starting where the old last column used to be,
# pushing the existing last column's data out to the new final column.
$ws.Range($lastColLetter + "1:" + $insertEndLetter + "1").EntireColumn.Insert()

# The old data (originally in $lastColLetter) now lives in the new final
# column; copy it back across the 4 freshly-inserted columns so every
# round column from $lastColLetter through the new final column repeats
# the same values.
$srcRange = $ws.Range($lastNewColLetter + "1:" + $lastNewColLetter + $lastRow)
$srcRange.Copy()

for ($i = 0; $i -lt $numNewCols; $i++) {
    $colLetter = Get-ColumnLetter ($lastCol + $i)
    $destRange = $ws.Range($colLetter + "1:" + $colLetter + $lastRow)
    $destRange.PasteSpecial()
}

$excel.CutCopyMode = 0
